$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the values that need to slide into new positions before they get overwritten.
$valB18 = $ws.Range("B18").Value
$valC18 = $ws.Range("C18").Value
$valB19 = $ws.Range("B19").Value
$valC19 = $ws.Range("C19").Value
$valB20 = $ws.Range("B20").Value
$valC20 = $ws.Range("C20").Value

# Remove row 13 entirely - everything below shifts up by one row.
$ws.Rows("13").Delete() | Out-Null

$daniela = "6666306 - Daniela Helena Pelegrine Guimarães"

# Objetivos: value replaced with the docente text.
$ws.Range("B10").Value = $daniela
$ws.Range("C10").Value = $daniela

# Programa resumido: value replaced with "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Programa: value replaced with the activation date.
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

# Método: value replaced with the docente text.
$ws.Range("B18").Value = $daniela
$ws.Range("C18").Value = $daniela

# Critério / Norma de recuperação / Bibliografia values shift up by one row.
$ws.Range("B19").Value = $valB18
$ws.Range("C19").Value = $valC18

$ws.Range("B20").Value = $valB19
$ws.Range("C20").Value = $valC19

$ws.Range("B21").Value = $valB20
$ws.Range("C21").Value = $valC20
